$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 466.46155
$ws.Cells.Item(33, 9).Value = 465.33334
$ws.Cells.Item(33, 11).Value = 465.33334
$ws.Cells.Item(33, 13).Value = -236.33334
$ws.Cells.Item(76, 8).Value = 5940.909
$ws.Cells.Item(76, 9).Value = 5166.8335
$ws.Cells.Item(76, 10).Value = 6231.1875
$ws.Cells.Item(76, 11).Value = 5166.8335
$ws.Cells.Item(76, 12).Value = 6231.1875
$ws.Cells.Item(76, 13).Value = -4851.8335
$ws.Cells.Item(76, 14).Value = -6861.1875
$ws.Cells.Item(79, 8).Value = 5940.909
$ws.Cells.Item(79, 9).Value = 5166.8335
$ws.Cells.Item(79, 10).Value = 6231.1875
$ws.Cells.Item(79, 11).Value = 5166.8335
$ws.Cells.Item(79, 12).Value = 6231.1875
$ws.Cells.Item(79, 13).Value = -4074.8335
$ws.Cells.Item(79, 14).Value = -8415.1875
$ws.Cells.Item(113, 8).Value = 2985.1667
$ws.Cells.Item(113, 10).Value = 3226.5
$ws.Cells.Item(113, 12).Value = 3226.5
$ws.Cells.Item(113, 14).Value = -9734.5
$ws.Cells.Item(115, 8).Value = 666.6667
$ws.Cells.Item(115, 9).Value = 666.6667
$ws.Cells.Item(115, 11).Value = 2000.0001
$ws.Cells.Item(115, 13).Value = -433.0001
$ws.Cells.Item(116, 8).Value = 3616.32
$ws.Cells.Item(116, 9).Value = 3367.0667
$ws.Cells.Item(116, 10).Value = 3990.2
$ws.Cells.Item(116, 11).Value = 3367.0667
$ws.Cells.Item(116, 12).Value = 3990.2
$ws.Cells.Item(116, 13).Value = 74.93330000000014
$ws.Cells.Item(116, 14).Value = -10874.2
$ws.Cells.Item(138, 8).Value = 1899.4933
$ws.Cells.Item(138, 9).Value = 1478.4286
$ws.Cells.Item(138, 10).Value = 1996.1311
$ws.Cells.Item(138, 11).Value = 4435.2858
$ws.Cells.Item(138, 12).Value = 5988.3933
$ws.Cells.Item(138, 13).Value = 704.7142000000003
$ws.Cells.Item(138, 14).Value = -16268.3933
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 1110.125
$ws.Cells.Item(45, 9).Value = 1123.5
$ws.Cells.Item(45, 10).Value = 1070
$ws.Cells.Item(45, 11).Value = 1123.5
$ws.Cells.Item(45, 12).Value = 1070
$ws.Cells.Item(45, 13).Value = -746.5
$ws.Cells.Item(45, 14).Value = -1824
$ws.Cells.Item(92, 8).Value = 2513944
$ws.Cells.Item(92, 10).Value = 2513944
$ws.Cells.Item(92, 12).Value = 2513944
$ws.Cells.Item(92, 14).Value = -2518936
$ws.Cells.Item(110, 8).Value = 1291.2858
$ws.Cells.Item(110, 9).Value = 918.1875
$ws.Cells.Item(110, 10).Value = 2485.2
$ws.Cells.Item(110, 11).Value = 918.1875
$ws.Cells.Item(110, 12).Value = 2485.2
$ws.Cells.Item(110, 13).Value = 1126.8125
$ws.Cells.Item(110, 14).Value = -6575.2
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 649.8
$ws.Cells.Item(80, 9).Value = 423.54544
$ws.Cells.Item(80, 10).Value = 926.3333
$ws.Cells.Item(80, 11).Value = 423.54544
$ws.Cells.Item(80, 12).Value = 926.3333
$ws.Cells.Item(80, 13).Value = 574.45456
$ws.Cells.Item(80, 14).Value = -2922.3333
$ws.Cells.Item(83, 8).Value = 649.8
$ws.Cells.Item(83, 9).Value = 423.54544
$ws.Cells.Item(83, 10).Value = 926.3333
$ws.Cells.Item(83, 11).Value = 2117.7272
$ws.Cells.Item(83, 12).Value = 4631.6665
$ws.Cells.Item(83, 13).Value = 2874.2728
$ws.Cells.Item(83, 14).Value = -14615.6665
$ws.Cells.Item(94, 8).Value = 35715576
$ws.Cells.Item(94, 9).Value = 41667668
$ws.Cells.Item(94, 11).Value = 41667668
$ws.Cells.Item(94, 13).Value = -41667217
$ws.Cells.Item(105, 8).Value = 83336104
$ws.Cells.Item(105, 9).Value = 90911840
$ws.Cells.Item(105, 11).Value = 90911840
$ws.Cells.Item(105, 13).Value = -90910093
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 62501090
$ws.Cells.Item(16, 9).Value = 76924090
$ws.Cells.Item(16, 10).Value = 1432.6666
$ws.Cells.Item(16, 11).Value = 76924090
$ws.Cells.Item(16, 12).Value = 1432.6666
$ws.Cells.Item(16, 13).Value = -76923803
$ws.Cells.Item(16, 14).Value = -2006.6666
$ws.Cells.Item(31, 8).Value = 1250.8298
$ws.Cells.Item(31, 9).Value = 1228.6444
$ws.Cells.Item(31, 11).Value = 1228.6444
$ws.Cells.Item(31, 13).Value = -933.6443999999999
$ws.Cells.Item(34, 8).Value = 1250.8298
$ws.Cells.Item(34, 9).Value = 1228.6444
$ws.Cells.Item(34, 11).Value = 1228.6444
$ws.Cells.Item(34, 13).Value = -1026.6444
$ws.Cells.Item(113, 8).Value = 62501090
$ws.Cells.Item(113, 9).Value = 76924090
$ws.Cells.Item(113, 10).Value = 1432.6666
$ws.Cells.Item(113, 11).Value = 76924090
$ws.Cells.Item(113, 12).Value = 1432.6666
$ws.Cells.Item(113, 13).Value = -76921920
$ws.Cells.Item(113, 14).Value = -5772.6666
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 705.64
$ws.Cells.Item(5, 9).Value = 618.6667
$ws.Cells.Item(5, 10).Value = 929.2857
$ws.Cells.Item(5, 11).Value = 1856.0001
$ws.Cells.Item(5, 12).Value = 2787.8571
$ws.Cells.Item(5, 13).Value = -1744.0001
$ws.Cells.Item(5, 14).Value = -3011.8571
$ws.Cells.Item(6, 8).Value = 576.25
$ws.Cells.Item(6, 9).Value = 270
$ws.Cells.Item(6, 10).Value = 1495
$ws.Cells.Item(6, 11).Value = 810
$ws.Cells.Item(6, 13).Value = -697
$ws.Cells.Item(6, 14).Value = -4711
$ws.Cells.Item(46, 8).Value = 1635.375
$ws.Cells.Item(46, 9).Value = 1726.1428
$ws.Cells.Item(46, 10).Value = 1000
$ws.Cells.Item(46, 11).Value = 5178.428400000001
$ws.Cells.Item(46, 14).Value = -3182
$ws.Cells.Item(68, 8).Value = 1420
$ws.Cells.Item(68, 9).Value = 1300
$ws.Cells.Item(68, 10).Value = 1450
$ws.Cells.Item(68, 11).Value = 3900
$ws.Cells.Item(68, 12).Value = 4350
$ws.Cells.Item(68, 13).Value = -3089
$ws.Cells.Item(68, 14).Value = -5972
$ws.Cells.Item(71, 8).Value = 1420
$ws.Cells.Item(71, 9).Value = 1300
$ws.Cells.Item(71, 10).Value = 1450
$ws.Cells.Item(71, 11).Value = 11700
$ws.Cells.Item(71, 12).Value = 13050
$ws.Cells.Item(71, 13).Value = -7644
$ws.Cells.Item(71, 14).Value = -21162
$ws.Cells.Item(92, 8).Value = 756.5714
$ws.Cells.Item(92, 9).Value = 839.2
$ws.Cells.Item(92, 10).Value = 550
$ws.Cells.Item(92, 11).Value = 2517.6
$ws.Cells.Item(92, 12).Value = 1650
$ws.Cells.Item(92, 13).Value = -1269.6
$ws.Cells.Item(92, 14).Value = -4146
$ws.Cells.Item(122, 8).Value = 702.5
$ws.Cells.Item(122, 9).Value = 420
$ws.Cells.Item(122, 10).Value = 823.5714
$ws.Cells.Item(122, 11).Value = 3780
$ws.Cells.Item(122, 12).Value = 7412.1426
$ws.Cells.Item(122, 13).Value = -1330
$ws.Cells.Item(122, 14).Value = -12312.1426
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(135, 8).Value = 705.64
$ws.Cells.Item(135, 9).Value = 618.6667
$ws.Cells.Item(135, 10).Value = 929.2857
$ws.Cells.Item(135, 11).Value = 5568.0003
$ws.Cells.Item(135, 12).Value = 8363.5713
$ws.Cells.Item(135, 13).Value = -3033.0003
$ws.Cells.Item(135, 14).Value = -13433.5713
$ws.Cells.Item(137, 8).Value = 2136.7896
$ws.Cells.Item(137, 9).Value = 981.8182
$ws.Cells.Item(137, 10).Value = 3724.875
$ws.Cells.Item(137, 11).Value = 2945.4546
$ws.Cells.Item(137, 12).Value = 11174.625
$ws.Cells.Item(137, 13).Value = 2154.5454
$ws.Cells.Item(137, 14).Value = -21374.625
$ws.Cells.Item(46, 8).Value = 14019.6
$ws.Cells.Item(46, 10).Value = 21433
$ws.Cells.Item(46, 12).Value = 21433
$ws.Cells.Item(46, 14).Value = -21745
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 859.5
$ws.Cells.Item(22, 9).Value = 538.3333
$ws.Cells.Item(22, 10).Value = 1100.375
$ws.Cells.Item(22, 11).Value = 538.3333
$ws.Cells.Item(22, 12).Value = 1100.375
$ws.Cells.Item(22, 13).Value = -243.3333
$ws.Cells.Item(22, 14).Value = -1690.375
$ws.Cells.Item(27, 8).Value = 859.5
$ws.Cells.Item(27, 9).Value = 538.3333
$ws.Cells.Item(27, 10).Value = 1100.375
$ws.Cells.Item(27, 11).Value = 538.3333
$ws.Cells.Item(27, 12).Value = 1100.375
$ws.Cells.Item(27, 13).Value = -431.3333
$ws.Cells.Item(27, 14).Value = -1314.375
$ws.Cells.Item(136, 8).Value = 1309.579
$ws.Cells.Item(136, 9).Value = 1323.5625
$ws.Cells.Item(136, 11).Value = 3970.6875
$ws.Cells.Item(136, 13).Value = -1420.6875
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(133, 8).Value = 45966.668
$ws.Cells.Item(133, 10).Value = 45966.668
$ws.Cells.Item(133, 12).Value = 45966.668
$ws.Cells.Item(133, 14).Value = -56086.668
